$wb = $excel.ActiveWorkbook

# --- Remove the empty "Sheet2" tab ---
$excel.DisplayAlerts = $false
$null = $wb.Worksheets.Item("Sheet2").Delete()

# --- Work on the "Program" sheet ---
$ws = $wb.Worksheets.Item("Program")

# Header row (A1:C1) becomes bold -> new font/style
$ws.Range("A1:C1").Font.Bold = $true

# New data rows (values entered in the same order the original
# author's shared-string table was built, so indices line up)
$ws.Range("A4").Value = 123
$ws.Range("B4").Value = 1345

$ws.Range("C5").Value = "Active"
$ws.Range("B8").Value = "4232#"
$ws.Range("B5").Value = "ZSY123 Desc"
$ws.Range("A5").Value = "CZX123"

$ws.Range("A6").Value = "DAZXS"
$ws.Range("B7").Value = "ABC Description123"
$ws.Range("A8").Value = 890364

# Selection / active cell ends on A8
$null = $ws.Range("A8").Select()

# Make "Program" the active tab
$null = $ws.Activate()
